# Insert two new price records right after the existing row 1238 (before the
# old row 1239), shifting every subsequent row down by two. This mirrors the
# weekly refresh: the two newest "Primera"/"Segunda" price points are added
# at the top of the historical block, pushing the rest of the series down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 1239:1301 down to 1241:1303, leaving two blank rows (1239, 1240)
# ready to receive the new data.
$ws.Rows.Item(1239).Resize(2).Insert()

# --- Row 1239 ("Primera") -------------------------------------------------
$ws.Cells.Item(1239, 1).Value2  = 3
$ws.Cells.Item(1239, 2).Value2  = "Femacal de La Calera"
$ws.Cells.Item(1239, 3).Value2  = "Coquimbo"
$ws.Cells.Item(1239, 4).Value2  = 45267
$ws.Cells.Item(1239, 5).Value2  = 5
$ws.Cells.Item(1239, 6).Value2  = 100112023
$ws.Cells.Item(1239, 7).Value2  = "Brócoli"
$ws.Cells.Item(1239, 8).Value2  = "Sin especificar"
$ws.Cells.Item(1239, 9).Value2  = "Primera"
$ws.Cells.Item(1239, 10).Value2 = 3100
$ws.Cells.Item(1239, 11).Value2 = 800
$ws.Cells.Item(1239, 12).Value2 = 900
$ws.Cells.Item(1239, 13).Value2 = 852
$ws.Cells.Item(1239, 14).Value2 = "$/unidad"
$ws.Cells.Item(1239, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(1239, 16).Value2 = 852
$ws.Cells.Item(1239, 17).Value2 = 1
$ws.Cells.Item(1239, 18).Value2 = "Hortaliza"

# --- Row 1240 ("Segunda") --------------------------------------------------
$ws.Cells.Item(1240, 1).Value2  = 3
$ws.Cells.Item(1240, 2).Value2  = "Femacal de La Calera"
$ws.Cells.Item(1240, 3).Value2  = "Coquimbo"
$ws.Cells.Item(1240, 4).Value2  = 45267
$ws.Cells.Item(1240, 5).Value2  = 5
$ws.Cells.Item(1240, 6).Value2  = 100112023
$ws.Cells.Item(1240, 7).Value2  = "Brócoli"
$ws.Cells.Item(1240, 8).Value2  = "Sin especificar"
$ws.Cells.Item(1240, 9).Value2  = "Segunda"
$ws.Cells.Item(1240, 10).Value2 = 1200
$ws.Cells.Item(1240, 11).Value2 = 700
$ws.Cells.Item(1240, 12).Value2 = 700
$ws.Cells.Item(1240, 13).Value2 = 700
$ws.Cells.Item(1240, 14).Value2 = "$/unidad"
$ws.Cells.Item(1240, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(1240, 16).Value2 = 700
$ws.Cells.Item(1240, 17).Value2 = 1
$ws.Cells.Item(1240, 18).Value2 = "Hortaliza"
